$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(9, 1).Value = 8995.64
$ws.Cells.Item(9, 2).Value = 8760
$ws.Cells.Item(9, 3).Value = 19.36
$ws.Cells.Item(9, 4).Value = 18.84
$ws.Cells.Item(9, 5).Value = $true
$ws.Cells.Item(9, 6).Value = -2.69
$ws.Cells.Item(9, 7).Value = 42612.672997685186
$ws.Cells.Item(9, 7).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(9, 8).Value = $true
